# Apply the "Integrating EU data files" changes to the Cargo Dist Conversion
# Factors workbook: replace the U.S.-centric notes/units with EU ones on the
# About sheet, and switch the CDCF-PMpPDOU / CDCF-FTMpFDOU conversion factors
# from "trillion" (10^12) units to billion (10^9) km-based units.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Remove old "For the U.S. model..." / unit rows (old rows 11-13) content
# and rebuild rows 11-17 with the new EU content + helper formula.
# NOTE: shared-string indices are assigned in order of first use, so the
# cells below are populated in the same order the reference workbook's
# shared string table expects (Freight before Passenger).
$about.Range("A11").Value = "For the EU model, the desired output units are:"
$about.Range("A13").Value = "Freight transport (Gtkm - Giga tonne-kilometers) - 10^9 tonne-kilometers - billion tonne-kilometers"
$about.Range("A12").Value = "Passenger transport (Gpkm - Giga passenger-kilometers) - 10^9 passenger-kilometers - billion passenger-kilometers"
$about.Range("A14").Value = ""
$about.Range("A15").Value = "ton-mile in the USA: 1 ton-mile * ( 0.907185 t / short ton) * ( 1.609344 km / mile ) = 1.460 tkm"
$about.Range("A16").Value = ""

$about.Range("A17").Formula = "=0.907185*1.609344"
$about.Range("A17").NumberFormat = "0.00000000000"

# ---------------------------------------------------------------------
# Sheet "CDCF-PMpPDOU"
# ---------------------------------------------------------------------
$pm = $wb.Worksheets.Item("CDCF-PMpPDOU")
$pm.Range("B2").Formula = "=1.609344*10^9"

# ---------------------------------------------------------------------
# Sheet "CDCF-FTMpFDOU"
# ---------------------------------------------------------------------
$ftm = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ftm.Range("B2").Formula = "=1.45997273664*10^9"
$ftm.Range("B5").NumberFormat = "0.00000000000"

$wb.Save()
